{"js": "// Replace each arithmetic-expression cell in the single table with its\n// updated equation (same layout, only the w:t text changes). The table is\n// 20 rows x 5 columns; we overwrite every cell's value in document order.\nconst newValues = [\n  [\"78-7=71\", \"72-49=23\", \"94-6=88\", \"32+6=38\", \"78-57=21\"],\n  [\"98-32=66\", \"77+8=85\", \"33+44=77\", \"68+18=86\", \"93+1=94\"],\n  [\"85-35=50\", \"11+78=89\", \"35+38=73\", \"72+3=75\", \"25+70=95\"],\n  [\"28+43=71\", \"83-81=2\", \"71-32=39\", \"31+4=35\", \"0+36=36\"],\n  [\"31+37=68\", \"58-16=42\", \"91-27=64\", \"58+11=69\", \"40+22=62\"],\n  [\"50-46=4\", \"50+35=85\", \"58+9=67\", \"84-4=80\", \"83+7=90\"],\n  [\"14+6=20\", \"59+27=86\", \"26-24=2\", \"44+21=65\", \"47-26=21\"],\n  [\"18+5=23\", \"69+7=76\", \"5+60=65\", \"92-56=36\", \"9+73=82\"],\n  [\"66-64=2\", \"88-10=78\", \"71-8=63\", \"30-14=16\", \"96+2=98\"],\n  [\"94-41=53\", \"35+26=61\", \"13+7=20\", \"17-13=4\", \"49-24=25\"],\n  [\"21+38=59\", \"38-11=27\", \"14+1=15\", \"10+89=99\", \"13+83=96\"],\n  [\"50+17=67\", \"89-59=30\", \"97-10=87\", \"12+0=12\", \"78-66=12\"],\n  [\"42-27=15\", \"72-22=50\", \"60-57=3\", \"99-89=10\", \"46-30=16\"],\n  [\"17-10=7\", \"86-5=81\", \"93-59=34\", \"68-42=26\", \"69-46=23\"],\n  [\"80-42=38\", \"8+77=85\", \"89-68=21\", \"60-55=5\", \"83-71=12\"],\n  [\"59-45=14\", \"28+5=33\", \"89-79=10\", \"60-25=35\", \"94-74=20\"],\n  [\"19+61=80\", \"16+3=19\", \"67-10=57\", \"15+84=99\", \"15+77=92\"],\n  [\"19+33=52\", \"23+66=89\", \"40+17=57\", \"77-38=39\", \"99-8=91\"],\n  [\"89-82=7\", \"94-2=92\", \"34-23=11\", \"30+29=59\", \"79-32=47\"],\n  [\"61+26=87\", \"99-23=76\", \"43+13=56\", \"27-0=27\", \"89-66=23\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\n// Sanity check: the grid shape must match the table before we overwrite it.\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Unexpected row count: table has ${table.rowCount}, expected ${newValues.length}`\n  );\n}\n\n// Writing the whole 2-D array in one shot updates each cell's text run\n// in place, leaving paragraph/run formatting (fonts, size, alignment)\n// untouched.\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace each arithmetic-expression cell in the single table with its\n# updated equation (same layout, only the cell text changes). The table\n# is 20 rows x 5 columns; $newValues[r][c] holds the replacement text for\n# row r+1 / column c+1 (Word COM collections are 1-indexed).\n$newValues = @(\n    @(\"78-7=71\", \"72-49=23\", \"94-6=88\", \"32+6=38\", \"78-57=21\"),\n    @(\"98-32=66\", \"77+8=85\", \"33+44=77\", \"68+18=86\", \"93+1=94\"),\n    @(\"85-35=50\", \"11+78=89\", \"35+38=73\", \"72+3=75\", \"25+70=95\"),\n    @(\"28+43=71\", \"83-81=2\", \"71-32=39\", \"31+4=35\", \"0+36=36\"),\n    @(\"31+37=68\", \"58-16=42\", \"91-27=64\", \"58+11=69\", \"40+22=62\"),\n    @(\"50-46=4\", \"50+35=85\", \"58+9=67\", \"84-4=80\", \"83+7=90\"),\n    @(\"14+6=20\", \"59+27=86\", \"26-24=2\", \"44+21=65\", \"47-26=21\"),\n    @(\"18+5=23\", \"69+7=76\", \"5+60=65\", \"92-56=36\", \"9+73=82\"),\n    @(\"66-64=2\", \"88-10=78\", \"71-8=63\", \"30-14=16\", \"96+2=98\"),\n    @(\"94-41=53\", \"35+26=61\", \"13+7=20\", \"17-13=4\", \"49-24=25\"),\n    @(\"21+38=59\", \"38-11=27\", \"14+1=15\", \"10+89=99\", \"13+83=96\"),\n    @(\"50+17=67\", \"89-59=30\", \"97-10=87\", \"12+0=12\", \"78-66=12\"),\n    @(\"42-27=15\", \"72-22=50\", \"60-57=3\", \"99-89=10\", \"46-30=16\"),\n    @(\"17-10=7\", \"86-5=81\", \"93-59=34\", \"68-42=26\", \"69-46=23\"),\n    @(\"80-42=38\", \"8+77=85\", \"89-68=21\", \"60-55=5\", \"83-71=12\"),\n    @(\"59-45=14\", \"28+5=33\", \"89-79=10\", \"60-25=35\", \"94-74=20\"),\n    @(\"19+61=80\", \"16+3=19\", \"67-10=57\", \"15+84=99\", \"15+77=92\"),\n    @(\"19+33=52\", \"23+66=89\", \"40+17=57\", \"77-38=39\", \"99-8=91\"),\n    @(\"89-82=7\", \"94-2=92\", \"34-23=11\", \"30+29=59\", \"79-32=47\"),\n    @(\"61+26=87\", \"99-23=76\", \"43+13=56\", \"27-0=27\", \"89-66=23\"),\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nif ($table.Rows.Count -ne $newValues.Count) {\n    throw \"Unexpected row count: table has $($table.Rows.Count), expected $($newValues.Count)\"\n}\n\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        # Cell.Range includes the trailing end-of-cell marker; assigning\n        # .Text replaces only the visible run text and keeps the run's\n        # existing character formatting (font/size).\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
